# Update profit files after running on 2025-11-03
# Appends the next day's row (row 78) to the profit log on Sheet1:
#   A78 = "11/03/2025"  (plain text, matching the existing date-as-text rows)
#   B78 = 9861.530000000001 (numeric profit value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date column to be stored as literal text (not auto-converted to
# a date serial number by Excel's input parsing), matching every other row
# in column A. Setting NumberFormat to Text ("@") before assigning the
# value prevents the "11/03/2025" -> date-serial auto-detection.
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = "11/03/2025"
# Drop the formatting we just applied so the new cell doesn't pick up an
# explicit style index (the source rows have no cell-level style either).
$ws.Range("A78").ClearFormats()

$ws.Range("B78").Value = 9861.530000000001
